$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("I2").Value = 6.2
$ws.Range("Q2").Value = 1.9
$ws.Range("U2").Value = 1.99
$ws.Range("AB2").Value = 8.4
# Row 3
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = 36
$ws.Range("H3").Value = 1.14
$ws.Range("I3").Value = 1.18
$ws.Range("J3").Value = 8.199999999999999
$ws.Range("N3").Value = 5.2
$ws.Range("P3").Value = 2.46
$ws.Range("Q3").Value = 1.53
$ws.Range("R3").Value = 1.58
$ws.Range("S3").Value = 2.32
$ws.Range("W3").Value = 1.03
$ws.Range("AB3").Value = 90
$ws.Range("AC3").Value = 26
$ws.Range("AD3").Value = 17
$ws.Range("AE3").Value = 19.5
$ws.Range("AH3").Value = 80
$ws.Range("AI3").Value = 85
# Row 4
$ws.Range("G4").Value = 3.05
$ws.Range("H4").Value = 2.58
$ws.Range("I4").Value = 3.55
$ws.Range("J4").Value = 2.8
$ws.Range("Q4").Value = 1.98
$ws.Range("V4").Value = 1.41
$ws.Range("W4").Value = 1.49
# Row 5
$ws.Range("F5").Value = 6.2
$ws.Range("G5").Value = 7.4
$ws.Range("H5").Value = 1.5
$ws.Range("I5").Value = 1.54
$ws.Range("J5").Value = 4.6
$ws.Range("L5").Value = 1.28
$ws.Range("N5").Value = 4.4
$ws.Range("P5").Value = 2.18
$ws.Range("Q5").Value = 1.71
$ws.Range("R5").Value = 1.45
$ws.Range("S5").Value = 2.78
$ws.Range("T5").Value = 1.83
$ws.Range("U5").Value = 1.96
$ws.Range("V5").Value = 2.84
$ws.Range("W5").Value = 1.16
$ws.Range("X5").Value = 21
$ws.Range("Z5").Value = 9.800000000000001
$ws.Range("AC5").Value = 11
$ws.Range("AF5").Value = 980
$ws.Range("AG5").Value = 27
$ws.Range("AH5").Value = 990
$ws.Range("AI5").Value = 34
$ws.Range("AJ5").Value = 250
$ws.Range("AK5").Value = 100
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 120
# Row 6
$ws.Range("G6").Value = 1.43
$ws.Range("H6").Value = 8.199999999999999
$ws.Range("I6").Value = 11
$ws.Range("J6").Value = 4.9
$ws.Range("T6").Value = 1.98
$ws.Range("U6").Value = 1.85
$ws.Range("V6").Value = 1.11
# Row 7
$ws.Range("T7").Value = 1.85
$ws.Range("AF7").Value = 11.5
# Row 8
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 10.5
$ws.Range("U8").Value = 1.78
# Row 9
$ws.Range("P9").Value = 1.78
# Row 10
$ws.Range("G10").Value = 4.1
$ws.Range("H10").Value = 2.1
$ws.Range("I10").Value = 2.22
$ws.Range("K10").Value = 3.65
$ws.Range("Q10").Value = 1.94
$ws.Range("S10").Value = 3.9
$ws.Range("V10").Value = 1.81
# Row 11
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 7.2
$ws.Range("I11").Value = 1.62
$ws.Range("J11").Value = 4.2
$ws.Range("Q11").Value = 2.14
$ws.Range("V11").Value = 2.62
# Row 12
$ws.Range("F12").Value = 2.88
$ws.Range("G12").Value = 2.9
$ws.Range("H12").Value = 2.86
$ws.Range("I12").Value = 2.88
$ws.Range("J12").Value = 3.25
$ws.Range("K12").Value = 3.3
$ws.Range("L12").Value = 1.53
$ws.Range("N12").Value = 2.94
$ws.Range("P12").Value = 1.64
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.24
$ws.Range("S12").Value = 5.1
$ws.Range("T12").Value = 2.06
$ws.Range("U12").Value = 1.89
$ws.Range("V12").Value = 1.53
$ws.Range("W12").Value = 1.52
$ws.Range("Y12").Value = 9
$ws.Range("AB12").Value = 9
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 60
$ws.Range("AJ12").Value = 44
$ws.Range("AM12").Value = 130
$ws.Range("AN12").Value = 44
# Row 13
$ws.Range("G13").Value = 4.6
$ws.Range("W13").Value = 1.28
# Row 14
$ws.Range("G14").Value = 3.55
$ws.Range("H14").Value = 2.44
$ws.Range("N14").Value = 2.68
$ws.Range("O14").Value = 1.53
$ws.Range("P14").Value = 1.58
$ws.Range("W14").Value = 1.39
$ws.Range("AO14").Value = 980
# Row 15
$ws.Range("F15").Value = 1.98
$ws.Range("G15").Value = 2.12
$ws.Range("I15").Value = 4.2
$ws.Range("J15").Value = 3.65
$ws.Range("N15").Value = 4.2
$ws.Range("P15").Value = 2.12
$ws.Range("Q15").Value = 1.74
$ws.Range("R15").Value = 1.44
$ws.Range("S15").Value = 2.92
$ws.Range("T15").Value = 1.67
$ws.Range("U15").Value = 2.24
$ws.Range("V15").Value = 1.32
$ws.Range("W15").Value = 1.89
$ws.Range("Z15").Value = 34
$ws.Range("AC15").Value = 11
$ws.Range("AM15").Value = 100
$ws.Range("AN15").Value = 15.5
# Row 16
$ws.Range("F16").Value = 2.28
$ws.Range("H16").Value = 3.6
$ws.Range("R16").Value = 1.35
$ws.Range("T16").Value = 1.83
$ws.Range("X16").Value = 13
$ws.Range("AE16").Value = 42
$ws.Range("AK16").Value = 23
$ws.Range("AN16").Value = 19
# Row 17
$ws.Range("F17").Value = 1.82
$ws.Range("G17").Value = 1.83
$ws.Range("I17").Value = 5.5
$ws.Range("K17").Value = 3.8
$ws.Range("V17").Value = 1.22
$ws.Range("W17").Value = 2.2
$ws.Range("AO17").Value = 95
# Row 18
$ws.Range("F18").Value = 3.75
$ws.Range("G18").Value = 4
$ws.Range("J18").Value = 3.45
$ws.Range("O18").Value = 1.37
$ws.Range("P18").Value = 1.83
$ws.Range("Q18").Value = 2.12
$ws.Range("V18").Value = 1.83
$ws.Range("W18").Value = 1.34
$ws.Range("AB18").Value = 13.5
$ws.Range("AF18").Value = 27
$ws.Range("AM18").Value = 120
$ws.Range("AN18").Value = 1000
$ws.Range("AO18").Value = 19.5
# Row 19
$ws.Range("G19").Value = 4
$ws.Range("P19").Value = 2.16
$ws.Range("R19").Value = 1.45
$ws.Range("T19").Value = 1.73
$ws.Range("AB19").Value = 16
$ws.Range("AM19").Value = 75
# Row 20
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 9.4
$ws.Range("N20").Value = 7
$ws.Range("T20").Value = 2.28
$ws.Range("U20").Value = 1.68
$ws.Range("Z20").Value = 280
$ws.Range("AB20").Value = 12.5
$ws.Range("AE20").Value = 460
$ws.Range("AI20").Value = 300
